$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header: doPhanGiaCameraSau -> doPhanGiaiCameraSau
$ws.Range("R1").Value = "doPhanGiaiCameraSau"

# Update model name: iPhone 30Cường -> iPhone 30
$ws.Range("A2").Value = "iPhone 30"

# Update camera type value: Plus -> Thường
$ws.Range("P2").Value = "Thường"

# Update the selection / view state
$ws.Range("H7").Select()
